# TASK_32, completed, Robot create game
# Adds the new "5.x - robot" task rows (35-39), extends the blank
# trailer rows (40-54), fixes up number formats on a few date cells,
# and moves the active-cell selection to D36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFmt = "dd/mm/yy\ h:mm;@"

# --- header row: give the date columns the same date format as column C ---
$ws.Range("D1:E1").NumberFormat = $dateFmt

# --- tidy up a couple of stray date formats in rows 32-33 ---
$ws.Range("E32").NumberFormat = $dateFmt
$ws.Range("D33:E33").NumberFormat = $dateFmt

# row 34 gains a bottom border (visually separates the old task list from
# the new "robot" task group being appended below it) - copy the border
# already used elsewhere on the sheet (row 6) so no new border style is
# invented
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Range("C6:E6").Copy() | Out-Null
$ws.Range("C34:E34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B34").Borders.Item(9).LineStyle = 1
$ws.Range("B34").Borders.Item(9).Weight = 2

# --- row 35: finish off TASK_32 ("5.1 - Создание игры с роботом") ---
$ws.Cells.Item(35, 2).Value = "5.1 – Создание игры с роботом. (LogicRobot.createGame)"
$ws.Cells.Item(35, 3).Value = 42013.770138888889
$ws.Cells.Item(35, 4).Value = 42013.770138888889
$ws.Cells.Item(35, 5).Value = 42013.80972222222
$ws.Range("C35:E35").NumberFormat = $dateFmt

# --- new rows 36-39: remaining robot tasks 5.2 - 5.5 ---
$ws.Cells.Item(36, 2).Value = "5.2 – Ход роботом. (LogicRobot.AIDoMove#1)"
$ws.Cells.Item(36, 3).Value = 42013.770138888889
$ws.Cells.Item(36, 4).Value = 42013.80972222222
$ws.Range("C36:D36").NumberFormat = $dateFmt
$ws.Cells.Item(36, 5).NumberFormat = $dateFmt

$ws.Cells.Item(37, 2).Value = "5.3 – Мозг робота. (LogicRobot.AIDoMove#2)"
$ws.Cells.Item(37, 3).Value = 42013.770138888889
$ws.Cells.Item(37, 3).NumberFormat = $dateFmt
$ws.Range("D37:E37").NumberFormat = $dateFmt

$ws.Cells.Item(38, 2).Value = "5.4 – Проверка победителя. (LogicRobot.checkWinner)"
$ws.Cells.Item(38, 3).Value = 42013.770138888889
$ws.Cells.Item(38, 3).NumberFormat = $dateFmt
$ws.Range("D38:E38").NumberFormat = $dateFmt

$ws.Cells.Item(39, 2).Value = "5.5 – Покинуть и игру с роботом. (LogicRobot.closeGame)"
$ws.Cells.Item(39, 3).Value = 42013.770138888889
$ws.Cells.Item(39, 3).NumberFormat = $dateFmt
$ws.Range("D39:E39").NumberFormat = $dateFmt

# column A running counter: extend the "+1" fill down through row 39
$ws.Cells.Item(37, 1).Formula = "=A36+1"
$ws.Cells.Item(38, 1).Formula = "=A37+1"
$ws.Cells.Item(39, 1).Formula = "=A38+1"

# --- rows 40-54: blank trailing rows, just carrying the date format ---
$ws.Range("D40:E54").NumberFormat = $dateFmt

# --- move the selection the author left it at ---
$ws.Range("D36").Select() | Out-Null
